$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "39.946.11"
$ws.Cells.Item(2, 5).Value = "  -0.27%  "

Set-TextCell 3 4 "2.217.74"
$ws.Cells.Item(3, 5).Value = "  -0.31%  "

$ws.Cells.Item(4, 5).Value = "  -0.07%  "

Set-TextCell 5 4 "293.38"
$ws.Cells.Item(5, 5).Value = "  +0.37%  "

Set-TextCell 6 4 "86.58"
$ws.Cells.Item(6, 5).Value = "  -2.03%  "

Set-TextCell 7 4 "0.511"
$ws.Cells.Item(7, 5).Value = "  -0.97%  "

$ws.Cells.Item(8, 5).Value = "  -0.07%  "

Set-TextCell 9 4 "0.466"
$ws.Cells.Item(9, 5).Value = "  -1.55%  "

Set-TextCell 10 4 "30.45"
$ws.Cells.Item(10, 5).Value = "  -1.93%  "

Set-TextCell 11 4 "0.0780"
$ws.Cells.Item(11, 5).Value = "  -0.99%  "

Set-TextCell 12 4 "50.20"
$ws.Cells.Item(12, 5).Value = "  +5.28%  "

$ws.Cells.Item(13, 5).Value = "  +3.26%  "

Set-TextCell 14 4 "6.41"
$ws.Cells.Item(14, 5).Value = "  -0.79%  "

Set-TextCell 15 4 "2.558.91"
$ws.Cells.Item(15, 5).Value = "  -0.44%  "

Set-TextCell 16 4 "13.79"
$ws.Cells.Item(16, 5).Value = "  -2.06%  "

Set-TextCell 17 4 "2.229.33"
$ws.Cells.Item(17, 5).Value = "  -0.29%  "

Set-TextCell 18 4 "0.733"
$ws.Cells.Item(18, 5).Value = "  +0.12%  "

Set-TextCell 19 4 "39.861.63"
$ws.Cells.Item(19, 5).Value = "  -0.37%  "

Set-TextCell 20 4 "0.0₃0883"
$ws.Cells.Item(20, 5).Value = "  -0.63%  "

Set-TextCell 21 4 "11.21"
$ws.Cells.Item(21, 5).Value = "  -7.69%  "

Set-TextCell 22 4 "5.75"
$ws.Cells.Item(22, 5).Value = "  -1.66%  "

Set-TextCell 23 4 "65.60"
$ws.Cells.Item(23, 5).Value = "  -0.23%  "

Set-TextCell 24 4 "235.83"
$ws.Cells.Item(24, 5).Value = "  +0.30%  "

$ws.Cells.Item(25, 5).Value = "  +0.13%  "

Set-TextCell 26 4 "2.46"
$ws.Cells.Item(26, 5).Value = "  -0.69%  "

Set-TextCell 27 4 "1.82"
$ws.Cells.Item(27, 5).Value = "  -3.09%  "

$ws.Cells.Item(28, 5).Value = "  +5.97%  "

Set-TextCell 29 4 "23.01"
$ws.Cells.Item(29, 5).Value = "  +1.19%  "

Set-TextCell 30 4 "9.25"
$ws.Cells.Item(30, 5).Value = "  -0.48%  "

Set-TextCell 31 4 "157.91"
$ws.Cells.Item(31, 5).Value = "  +3.56%  "

Set-TextCell 32 4 "31.52"
$ws.Cells.Item(32, 5).Value = "  -3.12%  "

$ws.Cells.Item(33, 5).Value = "  -0.01%  "

$ws.Cells.Item(34, 5).Value = "  -0.92%  "

Set-TextCell 35 4 "3.02"
$ws.Cells.Item(35, 5).Value = "  +5.07%  "

Set-TextCell 36 4 "0.0709"
$ws.Cells.Item(36, 5).Value = "  -1.76%  "

Set-TextCell 37 4 "2.33"
$ws.Cells.Item(37, 5).Value = "  -2.24%  "

Set-TextCell 38 4 "0.112"
$ws.Cells.Item(38, 5).Value = "  +0.14%  "

Set-TextCell 39 4 "0.0989"
$ws.Cells.Item(39, 5).Value = "  -0.99%  "

Set-TextCell 40 4 "1.73"
$ws.Cells.Item(40, 5).Value = "  +0.55%  "

Set-TextCell 41 4 "15.36"
$ws.Cells.Item(41, 5).Value = "  -5.10%  "

Set-TextCell 42 4 "2.083.86"
$ws.Cells.Item(42, 5).Value = "  -0.93%  "

Set-TextCell 43 4 "3.69"
$ws.Cells.Item(43, 5).Value = "  -4.45%  "

$ws.Cells.Item(44, 5).Value = "  -0.70%  "

$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 45 4 "9.97"
$ws.Cells.Item(45, 5).Value = "  -1.10%  "

$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 46 4 "17.78"
$ws.Cells.Item(46, 5).Value = "  -2.11%  "

Set-TextCell 47 4 "1.97"
$ws.Cells.Item(47, 5).Value = "  -9.31%  "

$ws.Cells.Item(48, 5).Value = "  -0.19%  "

Set-TextCell 49 4 "2.430.18"
$ws.Cells.Item(49, 5).Value = "  -0.65%  "

$ws.Cells.Item(50, 5).Value = "  +2.50%  "

Set-TextCell 51 4 "1.46"
$ws.Cells.Item(51, 5).Value = "  -0.27%  "
